$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("航天发展", "航天发展", "平潭发展")
    3  = @("中国卫星", "中国卫星", "航天发展")
    4  = @("神剑股份", "平潭发展", "中国卫星")
    5  = @("拓维信息", "神剑股份", "雷科防务")
    6  = @("平潭发展", "拓维信息", "神剑股份")
    7  = @("中超控股", "金风科技", "东百集团")
    8  = @("金风科技", "海南发展", "航天机电")
    9  = @("雷科防务", "雷科防务", "金风科技")
    10 = @("锋龙股份", "白银有色", "拓维信息")
    11 = @("天际股份", "中超控股", "浙江世宝")
    12 = @("翠微股份", "江西铜业", "翠微股份")
    13 = @("东百集团", "翠微股份", "中超控股")
    14 = @("海南发展", "御银股份", "国晟科技")
    15 = @("航天机电", "锋龙股份", "天奇股份")
    16 = @("超捷股份", "航天机电", "天际股份")
    17 = @("天奇股份", "东百集团", "泰尔股份")
    18 = @("白银有色", "天际股份", "海南发展")
    19 = @("浙江世宝", "拉卡拉", "再升科技")
    20 = @("江西铜业", "五洲新春", "雪人集团")
    21 = @("通宇通讯", "通宇通讯", "通宇通讯")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
}
